$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the header for the new third column (new shared string "두번째 변경점")
$ws.Range("C1").Value = "두번째 변경점"

# Add the new data values in column C
$ws.Range("C2").Value = 123
$ws.Range("C3").Value = 123
$ws.Range("C4").Value = 123

# Size column C to fit its content (closest achievable width to the target 13.75)
$ws.Columns.Item(3).ColumnWidth = 13

# Update the selected cell to match the target view
$ws.Range("E9").Select()
